# Automatische test-sync: 2025-08-18 20:48:50
#
# Adds the new log row (row 7) to the "Logs" sheet, extends the
# conditional-formatting ranges to include that row, and updates the
# "Dashboard" summary table (category counts) to reflect the new data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append the new e-mail log entry as row 7
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Interne taak"
$logs.Range("B7").Value = "kwaliteit@testbedrijf123.nl"
$logs.Range("C7").Value = "Leg dit even neer bij Koen."
$logs.Range("D7").Value = "Onbekend"
$logs.Range("E7").Value = "Fout bij verwerken (forward_to_fallback() got an unexpected keyword argument 'fallback_email')"
$logs.Range("F7").Value = "2025-08-18 20:48:17"
$logs.Range("G7").Value = "Nee"
$logs.Range("H7").Value = "Ja"
$logs.Range("I7").Value = "Nee"
$logs.Range("J7").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges from row 6 to row 7
# ---------------------------------------------------------------------
$ranges = @("D2:D6", "G2:G6", "H2:H6", "I2:I6", "J2:J6")
$newRanges = @("D2:D7", "G2:G7", "H2:H7", "I2:I7", "J2:J7")

for ($idx = 0; $idx -lt $ranges.Length; $idx++) {
    $oldRange = $ranges[$idx]
    $newRange = $newRanges[$idx]
    $fcs = $logs.Range($oldRange).FormatConditions
    $target = $logs.Range($newRange)
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($target)
    }
}

# ---------------------------------------------------------------------
# 3. "Dashboard" sheet: refresh the category summary counts
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Onbekend"
$dashboard.Range("B2").Value = 3
$dashboard.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$dashboard.Range("B3").Value = 2
